$d = $word.ActiveDocument

# --- Paragraph 2: "Asegurarse que el message procesor procesa solamente una
#     vez un evento." -> "Arreglar el memento para que funcione con los
#     ComplexEventSourced" (also drop the yellow highlight) -------------------
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.MoveEnd(1, -1)
$r2.HighlightColorIndex = 0

$rA = $p2.Range
$rA.MoveEnd(1, -1)
$rA.Find.Execute("Asegurarse que el message ", $true, $false, $false, `
    $false, $false, $true, 1, $false, `
    "Arreglar el memento para que funcione con los ", 2) | Out-Null

$rB = $p2.Range
$rB.MoveEnd(1, -1)
$rB.Find.Execute("procesor", $true, $false, $false, `
    $false, $false, $true, 1, $false, "ComplexEventSourced", 2) | Out-Null

$rC = $p2.Range
$rC.MoveEnd(1, -1)
$rC.Find.Execute(" procesa solamente una vez un evento.", $true, `
    $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Remove paragraph 3 entirely: "Colocar el lock statico en un command
#     handler que sea de alta contención" (also removes the old _GoBack
#     bookmark that lived inside it) -----------------------------------------
$d.Paragraphs.Item(3).Range.Delete() | Out-Null

# --- Remove paragraph 3 (was 4): "Asegurarse que si el complex event
#     procesor recibió un mensaje muy temprano..." ---------------------------
$d.Paragraphs.Item(3).Range.Delete() | Out-Null

# Paragraph 3 is now "Limitar la cantidad de mensajes que procesa al mismo
# tiempo (hacer tests de performance)" and needs no changes.

# --- Remove paragraph 4 (was 6): "Colocar el lock a los hight contentious
#     handlers." --------------------------------------------------------------
$d.Paragraphs.Item(4).Range.Delete() | Out-Null

# Paragraph 4 is now "Arreglar el event store rebuilder." Append a _GoBack
# bookmark right after its last character (before the paragraph mark).
$p4 = $d.Paragraphs.Item(4)
$rEnd = $p4.Range
$rEnd.MoveEnd(1, -1)
$rEnd.Collapse(0)
$rEnd.InsertAfter("X")
$rEnd.Collapse(0)
$rEnd.MoveStart(1, -1)
$d.Bookmarks.Add("_GoBack", $rEnd)
$rEnd.Text = ""
